$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 06:16"

# 2. Swap Austria / Australia country names (rows 71-72) and refresh Australia's stats
$ws.Range("A71").Value = "Australia"
$ws.Range("A72").Value = "Austria"

# 3. Swap Montserrat / Islas Malvinas country names (rows 213-214)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# 4. Refresh daily numbers (Casos totales, Nuevos casos, Casos activos, Recuperados,
#    Casos criticos, Muertes hoy, Muertes) for the countries whose figures changed

# India (row 6)
$ws.Range("B6").Value = 2647663
$ws.Range("C6").Value = 347
$ws.Range("D6").Value = 1919842
$ws.Range("E6").Value = 676776

# Pakistan (row 18)
$ws.Range("B18").Value = 289215
$ws.Range("C18").Value = 498
$ws.Range("D18").Value = 269087
$ws.Range("E18").Value = 13953
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 6175

# Kazajistan (row 29)
$ws.Range("D29").Value = 82777
$ws.Range("E29").Value = 18987

# Honduras (row 51)
$ws.Range("B51").Value = 50502
$ws.Range("C51").Value = 523
$ws.Range("D51").Value = 7339
$ws.Range("E51").Value = 41588
$ws.Range("G51").Value = 8
$ws.Range("H51").Value = 1575

# Row 71 - now Australia (updated stats)
$ws.Range("B71").Value = 23558
$ws.Range("C71").Value = 270
$ws.Range("D71").Value = 14080
$ws.Range("E71").Value = 9057
$ws.Range("G71").Value = 25
$ws.Range("H71").Value = 421

# Row 72 - now Austria (previous Austria stats, carried over unchanged)
$ws.Range("B72").Value = 23370
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 20681
$ws.Range("E72").Value = 1961
$ws.Range("H72").Value = 728

# Tailandia (row 117)
$ws.Range("B117").Value = 3378
$ws.Range("C117").Value = 1
$ws.Range("E117").Value = 126

# Row 213 - now Islas Malvinas (previous Islas Malvinas stats)
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214 - now Montserrat (previous Montserrat stats)
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
